$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
    if (-not $result) {
        Write-Host "WARNING: replace failed for:" $old
    }
}

function InsertAfter-Anchor($anchor, $newText) {
    $r = $d.Content
    $found = $r.Find.Execute($anchor)
    if (-not $found) {
        Write-Host "WARNING: anchor not found:" $anchor
        return
    }
    $r.Collapse(0)
    $r.InsertAfter($newText)
}

# --- Title ---
Replace-Text "Disentangling the Enigma of Dark Energy" "Delving into the Realm of Chemistry: Unveiling the Wonders at the Molecular Level"

# --- Author name ---
Replace-Text "Penelope Williams" "Alice Patterson"

# --- Email (split across two runs: local-part and domain) ---
Replace-Text "pwilliams@amail" "apatterson@eduvate"
Replace-Text "com" "org"

# --- Body paragraph 1 (four sentences) ---
Replace-Text "For centuries, the universe has captivated humankind's imagination, with its vastness and mysteries beckoning us to unravel its secrets" "Chemistry serves as an intriguing discipline that delves into the fundamental principles and interactions underlying the properties, composition, and transformations of matter"

Replace-Text " Of all the perplexing enigmas, dark energy stands out as one of the most enigmatic" " It acts as a magic key that opens the gates to the secrets hidden within substances, and it holds the power to unravel the mysteries that dictate how they behave"

Replace-Text " This mysterious force, comprising roughly 68% of the universe's total energy density, drives the universe's accelerated expansion and challenges our understanding of fundamental physics" " With fascination, we delve into the microscopic world, where atoms, the building blocks of all substances, engage in captivating dances of interactions"

Replace-Text " The pursuit of unraveling the nature of dark energy has launched an intellectual quest, captivating astrophysicists, cosmologists, and physicists worldwide, as they strive to comprehend its role in the universe's fate and dynamics" " In our journey of discovery, we perceive the mesmerizing interplay of particles as they exchange energy and form new combinations, thus weaving the intricate tapestry of the material world"

# --- Body paragraph 2 (three sentences, plus two new ones inserted at the end) ---
Replace-Text "The observed discrepancy between the expansion rate of the universe and the amount of matter it contains hints at the existence of a force counteracting the gravitational pull of matter" "Chemistry not only enables us to decipher the mysteries of matter but also empowers us to harness its potential for the betterment of society"

Replace-Text " This accelerating expansion suggests that either our understanding of gravity is incomplete or a new component, dark energy, is influencing the universe's evolution" " Through its lens, we uncover innovative solutions to global challenges, such as the development of cleaner energy sources, the engineering of advanced materials, and the synthesis of life-saving pharmaceuticals"

Replace-Text " The precise nature of dark energy remains elusive, as it interacts only through gravity and exerts a negative pressure, leading to the expansion of the universe" " Chemistry serves as an indispensable tool, aiding us in safeguarding the environment, enhancing human health, and facilitating technological advancements that shape our modern world"

# Insert two new sentences (period + sentence) right after the sentence above,
# before the existing trailing "." run.
InsertAfter-Anchor "Chemistry serves as an indispensable tool, aiding us in safeguarding the environment, enhancing human health, and facilitating technological advancements that shape our modern world" ". As we unveil the intricacies of chemistry, we unlock the potential to forge a sustainable and prosperous future"

# --- Body paragraph 3 (three sentences, plus four new ones inserted at the end) ---
Replace-Text "Dark energy's discovery has profoundly impacted cosmology, challenging prevailing theories and prompting the development of new cosmological models" "The versatility of chemistry extends to its diverse applications across various fields"

Replace-Text " Its existence implies that the universe's ultimate fate may be a ceaseless expansion or a sudden and catastrophic end" " It plays a profound role in the medical realm, facilitating the development of effective medications and therapies"

Replace-Text " The resolution of this cosmic conundrum holds the key to understanding the properties of dark energy and its influence on the universe's structure and evolution" " In agriculture, chemistry contributes to developing more productive crop varieties and devising innovative pest management strategies"

InsertAfter-Anchor " In agriculture, chemistry contributes to developing more productive crop varieties and devising innovative pest management strategies" ". Furthermore, it finds its place in materials science, leading to the creation of advanced materials with tailored properties that serve a multitude of purposes"
InsertAfter-Anchor " Furthermore, it finds its place in materials science, leading to the creation of advanced materials with tailored properties that serve a multitude of purposes" ". Chemistry's impact is apparent in energy production, propelling the transition to sustainable and efficient energy sources"
InsertAfter-Anchor " Chemistry's impact is apparent in energy production, propelling the transition to sustainable and efficient energy sources" ". The footprints of chemistry are ubiquitous in our daily lives beyond these core areas; from the clothes we wear, to the food we consume, to the devices we utilize, its presence is pervasive, making it a field of endless exploration and discovery"

# --- Summary heading: add lastRenderedPageBreak right before "Summary" ---
$sr = $d.Content
$sr.Find.Execute("Summary") | Out-Null
$sr.Collapse(1)
$breakAnchor = $d.Range($sr.Start, $sr.Start)
$lastRenderedPageBreakXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="28"/></w:rPr><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$breakAnchor.InsertXML($lastRenderedPageBreakXml) | Out-Null

# --- Summary paragraph ---
Replace-Text "Dark energy, an enigmatic force permeating the universe, constitutes roughly 68% of its total energy density" "Chemistry stands as a testament to the wonders of the microscopic world, offering boundless opportunities to explore the intricacies of matter and its transformations"

Replace-Text " Its existence challenges our comprehension of gravity and cosmology, hinting at the incompleteness of our physical understanding" " Its influence transcends the boundaries of scientific inquiry; it impacts various areas of human endeavor, including medicine, energy, agriculture, and materials science"

Replace-Text " The search for uncovering the nature of dark energy drives astrophysicists and cosmologists to explore new theories and models, seeking to decipher its properties and its intricate role in shaping the universe's destiny" " By delving into the realm of chemistry, we cultivate critical thinking abilities, nurture our curiosity, and gain a deeper understanding of the world around us"

Replace-Text " The resolution of this cosmic mystery promises to reshape our comprehension of the cosmos and its ultimate fate" " It is this pursuit of knowledge, coupled with an inquisitive spirit, that enables us to unravel the secrets of matter and harness its potential for societal progress"

# --- Add trailing empty paragraph at the very end of the document ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
